$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.720.54'
$ws.Range('E2').Value = '  +0.47%  '
$ws.Range('D3').Value = '2.471.56'
$ws.Range('E3').Value = '  +0.16%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '318.95'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '93.14'
$ws.Range('E6').Value = '  +1.37%  '
$ws.Range('E7').Value = '  +0.56%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.517'
$ws.Range('E9').Value = '  +0.32%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0868'
$ws.Range('E10').Value = '  +9.52%  '
$ws.Range('B11').Value = 'Avalanche'
$ws.Range('C11').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '33.40'
$ws.Range('E11').Value = '  +2.93%  '
$ws.Range('E12').Value = '  +0.65%  '
$ws.Range('D13').Value = '2.852.67'
$ws.Range('E13').Value = '  +0.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.91'
$ws.Range('E14').Value = '  +0.81%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.69'
$ws.Range('E15').Value = '  -1.48%  '
$ws.Range('D16').Value = '2.477.28'
$ws.Range('E16').Value = '  +0.46%  '
$ws.Range('E17').Value = '  +2.35%  '
$ws.Range('D18').Value = '41.676.23'
$ws.Range('E18').Value = '  +0.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.45'
$ws.Range('E19').Value = '  -0.90%  '
$ws.Range('E20').Value = '  +0.46%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.05'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.27'
$ws.Range('E22').Value = '  +1.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '239.86'
$ws.Range('E23').Value = '  +1.38%  '
$ws.Range('E24').Value = '  +0.80%  '
$ws.Range('E25').Value = '  +1.71%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.71'
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('E28').Value = '  +1.03%  '
$ws.Range('E29').Value = '  +0.90%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.11'
$ws.Range('E30').Value = '  +2.21%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '159.34'
$ws.Range('E31').Value = '  +2.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.49'
$ws.Range('E32').Value = '  +0.89%  '
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('E34').Value = '  +0.18%  '
$ws.Range('E35').Value = '  +0.61%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.53'
$ws.Range('E36').Value = '  +2.16%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.86'
$ws.Range('E37').Value = '  +4.63%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.92'
$ws.Range('E38').Value = '  +1.43%  '
$ws.Range('E39').Value = '  +1.98%  '
$ws.Range('E40').Value = '  +0.80%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.57'
$ws.Range('E41').Value = '  +6.20%  '
$ws.Range('E42').Value = '  +0.28%  '
$ws.Range('D43').Value = '1.998.30'
$ws.Range('E43').Value = '  +2.91%  '
$ws.Range('E44').Value = '  +0.60%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.85'
$ws.Range('E45').Value = '  -0.05%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.99'
$ws.Range('E46').Value = '  +2.88%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.53'
$ws.Range('E47').Value = '  +5.25%  '
$ws.Range('D48').Value = '2.710.30'
$ws.Range('E48').Value = '  +0.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '97.66'
$ws.Range('E49').Value = '  +0.35%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '74.45'
$ws.Range('E50').Value = '  +4.06%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '67.01'
$ws.Range('E51').Value = '  +0.04%  '
